$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3912.6
$ws.Range("I98").Value = 3745.3076
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 3745.3076
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -2247.3076
$ws.Range("N98").Value = -7996
$ws.Range("H116").Value = 3049
$ws.Range("I116").Value = 2233.3333
$ws.Range("J116").Value = 3398.5715
$ws.Range("K116").Value = 2233.3333
$ws.Range("L116").Value = 3398.5715
$ws.Range("M116").Value = 1208.6667
$ws.Range("N116").Value = -10282.5715
$ws.Range("H122").Value = 3912.6
$ws.Range("I122").Value = 3745.3076
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11235.9228
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8785.9228
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 194469.73
$ws.Range("I132").Value = 259005.23
$ws.Range("J132").Value = 863.2308
$ws.Range("K132").Value = 777015.6900000001
$ws.Range("L132").Value = 2589.6924
$ws.Range("M132").Value = -774485.6900000001
$ws.Range("N132").Value = -7649.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2943596.5
$ws.Range("I2").Value = 4267.3335
$ws.Range("J2").Value = 4203309
$ws.Range("K2").Value = 4267.3335
$ws.Range("L2").Value = 4203309
$ws.Range("M2").Value = -4154.3335
$ws.Range("N2").Value = -4203535
$ws.Range("H61").Value = 2680
$ws.Range("I61").Value = 1614.2858
$ws.Range("J61").Value = 5166.6665
$ws.Range("K61").Value = 1614.2858
$ws.Range("L61").Value = 5166.6665
$ws.Range("M61").Value = -1402.2858
$ws.Range("N61").Value = -5590.6665
$ws.Range("H116").Value = 2943596.5
$ws.Range("I116").Value = 4267.3335
$ws.Range("J116").Value = 4203309
$ws.Range("K116").Value = 4267.3335
$ws.Range("L116").Value = 4203309
$ws.Range("M116").Value = -1973.3335
$ws.Range("N116").Value = -4207897
$ws.Range("H136").Value = 2680
$ws.Range("I136").Value = 1614.2858
$ws.Range("J136").Value = 5166.6665
$ws.Range("K136").Value = 4842.857400000001
$ws.Range("L136").Value = 15499.9995
$ws.Range("M136").Value = -2292.857400000001
$ws.Range("N136").Value = -20599.9995
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H3").Value = 2943596.5
$ws.Range("I3").Value = 4267.3335
$ws.Range("J3").Value = 4203309
$ws.Range("K3").Value = 4267.3335
$ws.Range("L3").Value = 4203309
$ws.Range("M3").Value = -4153.3335
$ws.Range("N3").Value = -4203537

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1791.0714
$ws.Range("I86").Value = 1701.6897
$ws.Range("J86").Value = 1990.4615
$ws.Range("K86").Value = 1701.6897
$ws.Range("L86").Value = 1990.4615
$ws.Range("M86").Value = -578.6896999999999
$ws.Range("N86").Value = -4236.461499999999
$ws.Range("H89").Value = 1791.0714
$ws.Range("I89").Value = 1701.6897
$ws.Range("J89").Value = 1990.4615
$ws.Range("K89").Value = 8508.448499999999
$ws.Range("L89").Value = 9952.307499999999
$ws.Range("M89").Value = -2892.448499999999
$ws.Range("N89").Value = -21184.3075

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2303.3333
$ws.Range("I56").Value = 2303.3333
$ws.Range("K56").Value = 2303.3333
$ws.Range("M56").Value = -1773.3333
$ws.Range("H96").Value = 141411410
$ws.Range("J96").Value = 141411410
$ws.Range("L96").Value = 424234230
$ws.Range("N96").Value = -424238348
$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 9000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -17180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5131.8945
$ws.Range("I70").Value = 5338.5
$ws.Range("J70").Value = 4981.636
$ws.Range("K70").Value = 5338.5
$ws.Range("L70").Value = 4981.636
$ws.Range("M70").Value = -5068.5
$ws.Range("N70").Value = -5521.636
$ws.Range("H73").Value = 5131.8945
$ws.Range("I73").Value = 5338.5
$ws.Range("J73").Value = 4981.636
$ws.Range("K73").Value = 5338.5
$ws.Range("L73").Value = 4981.636
$ws.Range("M73").Value = -4402.5
$ws.Range("N73").Value = -6853.636
$ws.Range("H88").Value = 92588.25
$ws.Range("I88").Value = 89177
$ws.Range("J88").Value = 95999.5
$ws.Range("K88").Value = 89177
$ws.Range("L88").Value = 95999.5
$ws.Range("M88").Value = -88726
$ws.Range("N88").Value = -96901.5
$ws.Range("H91").Value = 92588.25
$ws.Range("I91").Value = 89177
$ws.Range("J91").Value = 95999.5
$ws.Range("K91").Value = 89177
$ws.Range("L91").Value = 95999.5
$ws.Range("M91").Value = -87617
$ws.Range("N91").Value = -99119.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 3000
$ws.Range("J5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3226
$ws.Range("H22").Value = 391.42856
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 406.66666
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 406.66666
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -996.66666
$ws.Range("H27").Value = 391.42856
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 406.66666
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 406.66666
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -620.66666
$ws.Range("H31").Value = 1271.1428
$ws.Range("I31").Value = 632.6667
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 632.6667
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -384.6667
$ws.Range("N31").Value = -2246
$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15580
$ws.Range("H34").Value = 6474.5
$ws.Range("I34").Value = 950
$ws.Range("J34").Value = 11999
$ws.Range("K34").Value = 950
$ws.Range("L34").Value = 11999
$ws.Range("M34").Value = -778
$ws.Range("N34").Value = -12343
$ws.Range("H35").Value = 404.125
$ws.Range("I35").Value = 404.125
$ws.Range("K35").Value = 404.125
$ws.Range("M35").Value = -68.125
$ws.Range("H38").Value = 18000
$ws.Range("J38").Value = 18000
$ws.Range("L38").Value = 18000
$ws.Range("N38").Value = -18820
$ws.Range("H41").Value = 14999
$ws.Range("J41").Value = 14998
$ws.Range("L41").Value = 14998
$ws.Range("N41").Value = -15874
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
